$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-08-15"

# Update header label in I1
$ws.Range("I1").Value = "2022 (through 08-15)"

# Update monthly figures for the "2022 (through 08-15)" column (column I)
$ws.Range("I7").Value = 143   # June
$ws.Range("I9").Value = 85    # August

# Update the running total
$ws.Range("I14").Value = 1056
